$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update s6_cb (column H) and weighted_total (column N) for rows 2-5
# following the removal of a deleted question from the underlying data.

$ws.Range("H2").Value = 0.4444444444444444
$ws.Range("N2").Value = 0.3102380952380953

$ws.Range("H3").Value = 0.3333333333333333
$ws.Range("N3").Value = 0.2751190476190477

$ws.Range("H4").Value = 0.3333333333333333
$ws.Range("N4").Value = 0.2488095238095238

$ws.Range("H5").Value = 0.2222222222222222
$ws.Range("N5").Value = 0.1947619047619047
